# "Opdateret d. 02-12-2025" -> "Opdateret d. 05-12-2025"
#
# Weekly refresh of the approved-systems sheet: the sheet name carries
# the "last updated" date. Renaming the sheet keeps the workbook-scoped
# defined name "Lægesystemer___modtage" in sync automatically, since its
# formula is the sheet-qualified range
# 'Opdateret d. 02-12-2025'!$A$1:$K$40.
$wb = $excel.ActiveWorkbook

$oldName = "Opdateret d. 02-12-2025"
$newName = "Opdateret d. 05-12-2025"

try {
    $ws = $wb.Sheets.Item($oldName)
} catch {
    $ws = $wb.ActiveSheet
}

$ws.Name = $newName
